# Refactoring des CSV-Importers als Vorbereitung fuer den Import von Sprints:
# add a new "Sprint" column (F) to the PBL sheet and adjust the estimate
# for one PBI (D5) that was split up.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# New column F header + values (these are brand-new shared strings, so the
# order they are first written in matches the order they are appended to
# sharedStrings.xml: Sprint, Sprint 1, Sprint 1 + Sprint 2, Sprint 2).
$ws.Range("F1").Value = "Sprint"
$ws.Range("F2").Value = "Sprint 1"
$ws.Range("F3").Value = "Sprint 1, Sprint 2"
$ws.Range("F4").Value = "Sprint 2"
$ws.Range("F5").Value = "Sprint 2"

# Column F visual width (matches the bestFit width used for the other
# narrow columns in this sheet).
$ws.Columns.Item(6).ColumnWidth = 14.7

# PBI 6 (row 5) is now estimated at half a story point instead of 2.
$ws.Range("D5").Value = 0.5

# Selection left on D3:D5 (active cell D3) after the edit.
$ws.Range("D3:D5").Select()
